# The deck's theme ("Integral" colour scheme, stored in ppt/theme/theme1.xml
# and linked from the slide master) is swapped for the stock "Office Theme"
# colour scheme (the set of colours that, before this edit, lived unused in
# ppt/theme/theme2.xml, linked only from the notes master).
#
# PowerPoint's theme colours are exposed on the COM object model through
# Master.Theme.ThemeColorScheme, whose 12 items map 1:1 to the OOXML
# <a:clrScheme> children in document order:
#   1 dk1  2 lt1  3 dk2  4 lt2  5 accent1 6 accent2 7 accent3
#   8 accent4 9 accent5 10 accent6 11 hlink 12 folHlink
# ThemeColor.RGB uses the same packed 0xBBGGRR integer as the VBA RGB()
# function, so a small helper converts the target "RRGGBB" hex strings.

function Convert-HexToRgbInt($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r -bor ($g -shl 8) -bor ($b -shl 16)
}

# Target "Office Theme" colour scheme, in <a:clrScheme> order.
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

for ($i = 0; $i -lt $officeThemeColors.Count; $i++) {
    $colorScheme.Item($i + 1).RGB = Convert-HexToRgbInt $officeThemeColors[$i]
}
